$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.621.06'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '2.294.44'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.82'
$ws.Range("E5").Value = '  +18.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.15'
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.615'
$ws.Range("E9").Value = '  +0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.73'
$ws.Range("E10").Value = '  +4.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0939'
$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.58'
$ws.Range("E12").Value = '  +10.22%  '

$ws.Range("E13").Value = '  +1.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.60'
$ws.Range("E14").Value = '  +2.99%  '

$ws.Range("D15").Value = '2.633.92'
$ws.Range("E15").Value = '  -0.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.849'
$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").Value = '2.291.52'
$ws.Range("E17").Value = '  -0.19%  '

$ws.Range("D18").Value = '43.639.26'
$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("E19").Value = '  +2.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.57'
$ws.Range("E20").Value = '  +6.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.48'
$ws.Range("E21").Value = '  +1.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.51'
$ws.Range("E22").Value = '  +3.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.50'
$ws.Range("E23").Value = '  +0.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.60'
$ws.Range("E24").Value = '  +5.08%  '

$ws.Range("E25").Value = '  +14.39%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.54'
$ws.Range("E27").Value = '  +3.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.34'
$ws.Range("E28").Value = '  +5.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.41'
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '176.46'
$ws.Range("E31").Value = '  +0.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.65'
$ws.Range("E32").Value = '  -0.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0925'
$ws.Range("E33").Value = '  +4.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.54'
$ws.Range("E34").Value = '  +3.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  +0.81%  '

$ws.Range("E36").Value = '  +9.37%  '

$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0356'
$ws.Range("E38").Value = '  +0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.81'
$ws.Range("E39").Value = '  +12.78%  '

$ws.Range("B40").Value = 'MultiversX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '74.13'
$ws.Range("E40").Value = '  +15.53%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.42'
$ws.Range("E41").Value = '  +4.77%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.81'
$ws.Range("E42").Value = '  +12.94%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.243'
$ws.Range("E43").Value = '  +3.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.43'
$ws.Range("E44").Value = '  +6.47%  '

$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.97'
$ws.Range("E46").Value = '  +14.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.78'
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.21'
$ws.Range("E48").Value = '  +5.62%  '

$ws.Range("E49").Value = '  -1.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.24'
$ws.Range("E50").Value = '  +4.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.445'
$ws.Range("E51").Value = '  +4.42%  '

